$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the value 5 for the newly-scored cells in column E
$ws.Range("E10").Value = 5
$ws.Range("E12").Value = 5
$ws.Range("E15").Value = 5
$ws.Range("E17").Value = 5
$ws.Range("E23").Value = 5
$ws.Range("E25").Value = 5
$ws.Range("E26").Value = 5
$ws.Range("E27").Value = 5
$ws.Range("E28").Value = 5

# Move the active selection to E11 (the bottom-right frozen pane follows)
$ws.Range("E11").Select()
